$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 50522576
$ws.Range("D2").Value = 26543024
$ws.Range("E2").Value = 16274934
$ws.Range("F2").Value = 47.52
$ws.Range("G2").Value = 13545934
$ws.Range("H2").Value = 104.22
